$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 2797

$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 0

$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 127
$ws.Range("D77").Value = 11

$ws.Range("D107").Value = 0
$ws.Range("D108").Value = 0
$ws.Range("D109").Value = 0
$ws.Range("D110").Value = 0
$ws.Range("D111").Value = 0
$ws.Range("D112").Value = 0
$ws.Range("D113").Value = 0

$ws.Range("D129").Value = 0
$ws.Range("D130").Value = 0
